# Feb 12 Update 2
# Insert two new rows at the top of the survey list (rows 2-3) for the
# newest surveys (Survey 44 / Survey 43), shifting all existing survey
# rows down by two, then update the sheet view selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 2 (pushes old rows 2..43 -> 4..45)
$ws.Rows("2:3").Insert()

# New row 2: Survey 44 baseline, Feb 9 - Feb 11
$ws.Range("A2").Value = "Survey 44"
$ws.Range("A3").Value = "Survey 43"
$ws.Range("B3").Value = "Feb 2 - Feb 4"
$ws.Range("B2").Value = "Feb 9 - Feb 11"
$ws.Range("E3").Value = "Social Support`nHealthcare`nVaccines`nRISER`nUnemployment`nStimulus"
$ws.Range("C2").Value = "x"
$ws.Range("D3").Value = "x"

# Row 3 (Survey 43 follow-up row) needs extra height for the module note
$ws.Rows("3:3").RowHeight = 102

# The newly inserted row 2 picked up a stray formatted-but-empty cell in
# column E (inherited from the row below); clear it so it doesn't persist.
$ws.Range("E2").Clear()

# Update the saved view/selection state
$ws.Range("E5").Select()
